$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column E ("Quantidade restante") shifting the old PN/Nome
# columns to F/G. Excel copies the left neighbour's (D) formatting into the
# freshly inserted column, which already reproduces the boxed border look
# the diff expects for the new column.
$ws.Columns.Item(5).Insert()

# Header for the new column.
$ws.Range("E3").Value = "Quantidade restante"

# Row 4 gets a standalone formula; rows 5-13 are filled down from it so the
# engine stores them as one shared-formula group (matches how the workbook
# was actually edited in Excel: type the formula once, then fill down).
$ws.Range("E4").Formula = "=B4-D4"
$ws.Range("E5:E13").Formula = "=B5-D5"

# The whole card (title + table) got a bit wider; reapply a medium (thick)
# outside-left edge to column A so it again mirrors the medium edge already
# present on the table's right side.
$leftEdge = $ws.Range("A1:A13").Borders.Item(7)
$leftEdge.LineStyle = 1
$leftEdge.Weight = -4138

# Match the column width Excel would give the new column (same width as the
# column it was duplicated from).
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Restore the cursor position recorded in the saved file.
$ws.Range("I17").Select()
